$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has columns: A=code, B=status, C=codeforiati:group-code, D=codeforiati:group-name
# The edit swaps the order of the group-code / group-name pair (and corresponding
# per-group code/name values) within the shared-strings table, which in effect
# means columns C and D need to swap their header and values so that
# "group-name" (and the country/org name) comes before "group-code" (and the
# country/org code) in the string table, while the underlying data per row
# stays a {code, name} pair, just reordered as {name, code}.

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
